# Append: 2026-01-02 12:50 JST
# Update the "取得日時" (retrieved at) timestamp in column A for rows 2-8
# on the "ランサーズ" sheet from "2026-01-02 12:37:43" to "2026-01-02 12:50:07".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-02 12:50:07"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
